$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1408.3334
$ws.Range("I18").Value = 690
$ws.Range("J18").Value = 5000
$ws.Range("K18").Value = 690
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = -406
$ws.Range("N18").Value = -5568

$ws.Range("H106").Value = 5794.684
$ws.Range("I106").Value = 3220.1875
$ws.Range("J106").Value = 19525.334
$ws.Range("K106").Value = 3220.1875
$ws.Range("L106").Value = 19525.334
$ws.Range("M106").Value = -2589.1875
$ws.Range("N106").Value = -20787.334

$ws.Range("H107").Value = 803
$ws.Range("J107").Value = 1117
$ws.Range("L107").Value = 1117
$ws.Range("N107").Value = -4957

$ws.Range("H132").Value = 5807.641
$ws.Range("I132").Value = 5965.222
$ws.Range("K132").Value = 17895.666
$ws.Range("M132").Value = -15365.666

$ws.Range("H137").Value = 1830.0625
$ws.Range("I137").Value = 1243
$ws.Range("J137").Value = 2137.5715
$ws.Range("K137").Value = 3729
$ws.Range("L137").Value = 6412.7145
$ws.Range("M137").Value = -1179
$ws.Range("N137").Value = -11512.7145

$ws.Range("H138").Value = 3758.7693
$ws.Range("I138").Value = 3769.12
$ws.Range("K138").Value = 11307.36
$ws.Range("M138").Value = -6167.360000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 37820.168
$ws.Range("I32").Value = 39707.594
$ws.Range("K32").Value = 39707.594
$ws.Range("M32").Value = -39420.594

$ws.Range("H39").Value = 12665.182
$ws.Range("I39").Value = 13687.375
$ws.Range("J39").Value = 9939.333000000001
$ws.Range("K39").Value = 13687.375
$ws.Range("L39").Value = 9939.333000000001
$ws.Range("M39").Value = -13167.375
$ws.Range("N39").Value = -10979.333

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("N40").Value = 0

$ws.Range("H45").Value = 2426.9
$ws.Range("I45").Value = 1334.4546
$ws.Range("K45").Value = 1334.4546
$ws.Range("M45").Value = -957.4546

$ws.Range("H61").Value = 2349.6
$ws.Range("I61").Value = 2397.5
$ws.Range("K61").Value = 2397.5
$ws.Range("M61").Value = -2185.5

$ws.Range("H74").Value = 1745.4783
$ws.Range("I74").Value = 1221.5834
$ws.Range("J74").Value = 2317
$ws.Range("K74").Value = 1221.5834
$ws.Range("L74").Value = 2317
$ws.Range("M74").Value = -347.5834
$ws.Range("N74").Value = -4065

$ws.Range("H77").Value = 1745.4783
$ws.Range("I77").Value = 1221.5834
$ws.Range("J77").Value = 2317
$ws.Range("K77").Value = 6107.916999999999
$ws.Range("L77").Value = 11585
$ws.Range("M77").Value = -1739.916999999999
$ws.Range("N77").Value = -20321

$ws.Range("H97").Value = 7940.467
$ws.Range("J97").Value = 2208.25
$ws.Range("L97").Value = 2208.25
$ws.Range("N97").Value = -3200.25

$ws.Range("H110").Value = 3963.3333
$ws.Range("I110").Value = 3445
$ws.Range("K110").Value = 3445
$ws.Range("M110").Value = -1400

$ws.Range("H132").Value = 43631.707
$ws.Range("I132").Value = 45441.78
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 136325.34
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -133795.34
$ws.Range("N132").Value = -11060

$ws.Range("H136").Value = 2349.6
$ws.Range("I136").Value = 2397.5
$ws.Range("K136").Value = 7192.5
$ws.Range("M136").Value = -4642.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 92641.73
$ws.Range("J22").Value = 1189.5
$ws.Range("L22").Value = 1189.5
$ws.Range("N22").Value = -1535.5

$ws.Range("H86").Value = 4972.5835
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 4972.5835
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 4972.5835
$ws.Range("N86").Value = -7218.5835

$ws.Range("H89").Value = 4972.5835
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 4972.5835
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 24862.9175
$ws.Range("N89").Value = -36094.9175

$ws.Range("H107").Value = 3061.3333
$ws.Range("I107").Value = 1783.7142
$ws.Range("J107").Value = 4850
$ws.Range("K107").Value = 1783.7142
$ws.Range("L107").Value = 4850
$ws.Range("M107").Value = 136.2858000000001
$ws.Range("N107").Value = -8690

$ws.Range("H134").Value = 2410.1562
$ws.Range("I134").Value = 2262.0967
$ws.Range("K134").Value = 6786.2901
$ws.Range("M134").Value = -4251.2901

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 61072.117
$ws.Range("I58").Value = 64607
$ws.Range("K58").Value = 64607
$ws.Range("M58").Value = -64404

$ws.Range("H59").Value = 67367
$ws.Range("I59").Value = 2104
$ws.Range("K59").Value = 2104
$ws.Range("M59").Value = -959

$ws.Range("H62").Value = 4276
$ws.Range("J62").Value = 3006
$ws.Range("L62").Value = 3006
$ws.Range("N62").Value = -4254

$ws.Range("H65").Value = 4276
$ws.Range("J65").Value = 3006
$ws.Range("L65").Value = 15030
$ws.Range("N65").Value = -21270

$ws.Range("H68").Value = 55000
$ws.Range("J68").Value = 55000
$ws.Range("L68").Value = 55000
$ws.Range("N68").Value = -56498

$ws.Range("H71").Value = 55000
$ws.Range("J71").Value = 55000
$ws.Range("L71").Value = 165000
$ws.Range("N71").Value = -172488

$ws.Range("H103").Value = 59213.75
$ws.Range("I103").Value = 49261.5
$ws.Range("J103").Value = 69166
$ws.Range("K103").Value = 49261.5
$ws.Range("L103").Value = 69166
$ws.Range("M103").Value = -48089.5
$ws.Range("N103").Value = -71510

$ws.Range("H136").Value = 61072.117
$ws.Range("I136").Value = 64607
$ws.Range("K136").Value = 193821
$ws.Range("M136").Value = -191271

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 13216.5
$ws.Range("I70").Value = 10859.8
$ws.Range("J70").Value = 25000
$ws.Range("K70").Value = 32579.4
$ws.Range("L70").Value = 75000
$ws.Range("M70").Value = -32264.4
$ws.Range("N70").Value = -75630

$ws.Range("H73").Value = 13216.5
$ws.Range("I73").Value = 10859.8
$ws.Range("J73").Value = 25000
$ws.Range("K73").Value = 32579.4
$ws.Range("L73").Value = 75000
$ws.Range("M73").Value = -31487.4
$ws.Range("N73").Value = -77184

$ws.Range("H92").Value = 295.07144
$ws.Range("I92").Value = 248.22223
$ws.Range("J92").Value = 379.4
$ws.Range("K92").Value = 744.66669
$ws.Range("L92").Value = 1138.2
$ws.Range("M92").Value = 503.33331
$ws.Range("N92").Value = -3634.2

$ws.Range("H140").Value = 2428.9656
$ws.Range("I140").Value = 2101.6667
$ws.Range("K140").Value = 6305.000100000001
$ws.Range("M140").Value = -1125.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3897.9412
$ws.Range("I102").Value = 2646.8
$ws.Range("J102").Value = 5685.2856
$ws.Range("K102").Value = 2646.8
$ws.Range("L102").Value = 5685.2856
$ws.Range("M102").Value = -1024.8
$ws.Range("N102").Value = -8929.285599999999

$ws.Range("H107").Value = 77968.92
$ws.Range("J107").Value = 1433.2858
$ws.Range("L107").Value = 1433.2858
$ws.Range("N107").Value = -5273.2858

$ws.Range("H132").Value = 57650.445
$ws.Range("I132").Value = 60688.707
$ws.Range("K132").Value = 182066.121
$ws.Range("M132").Value = -179536.121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 363.75
$ws.Range("I9").Value = 363.75
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 363.75
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -139.75

$ws.Range("H93").Value = 2774.75
$ws.Range("I93").Value = 2200
$ws.Range("J93").Value = 3119.6
$ws.Range("K93").Value = 2200
$ws.Range("L93").Value = 3119.6
$ws.Range("M93").Value = -952
$ws.Range("N93").Value = -5615.6

$ws.Range("H132").Value = 85119.13
$ws.Range("I132").Value = 104082.336
$ws.Range("J132").Value = 9266.333000000001
$ws.Range("K132").Value = 312247.008
$ws.Range("L132").Value = 27798.999
$ws.Range("M132").Value = -309717.008
$ws.Range("N132").Value = -32858.999

$ws.Range("H136").Value = 2972.6428
$ws.Range("I136").Value = 1910
$ws.Range("K136").Value = 5730
$ws.Range("M136").Value = -3180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 3000000
$ws.Range("I29").Value = 3000000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 3000000
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -2999710

$ws.Range("H109").Value = 96606.28999999999
$ws.Range("J109").Value = 96606.28999999999
$ws.Range("L109").Value = 96606.28999999999
$ws.Range("N109").Value = -99380.28999999999

$ws.Range("H122").Value = 5749.6523
$ws.Range("I122").Value = 6212.15
$ws.Range("J122").Value = 2666.3333
$ws.Range("K122").Value = 18636.45
$ws.Range("L122").Value = 7998.999899999999
$ws.Range("M122").Value = -16186.45
$ws.Range("N122").Value = -12898.9999

$ws.Range("H132").Value = 75746.71000000001
$ws.Range("I132").Value = 75746.71000000001
$ws.Range("K132").Value = 227240.13
$ws.Range("M132").Value = -224710.13
